# dataloader modified and environment data only train/infer supported
#
# The "date" column (A) was previously stored as a real Excel date serial
# number formatted with a custom date/time numFmt (style index 2, which
# used numFmtId 165 "YYYY-MM-DD HH:MM:SS"). This edit converts those cells
# to plain numeric values written directly in YYYYMMDD form (e.g. 43227 /
# "2018-05-07" becomes the literal number 20180507) and removes the
# special date formatting/style from those cells so they fall back to the
# default (unstyled) cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blocks of rows in column A that share the same YYYYMMDD replacement value
# (each block is one weekly batch of 16 product rows).
$groups = @(
    @{ Range = "A2:A17";    Value = 20180507 },
    @{ Range = "A18:A33";   Value = 20180514 },
    @{ Range = "A34:A49";   Value = 20180521 },
    @{ Range = "A50:A65";   Value = 20180528 },
    @{ Range = "A66:A81";   Value = 20180604 },
    @{ Range = "A82:A97";   Value = 20180611 },
    @{ Range = "A98:A113";  Value = 20180618 },
    @{ Range = "A114:A129"; Value = 20180625 },
    @{ Range = "A130:A145"; Value = 20180702 }
)

foreach ($g in $groups) {
    $rng = $ws.Range($g.Range)
    # Reset to the default (unstyled) cell style before assigning the plain
    # numeric value, so the custom date numFmt (s="2", numFmtId 165) is no
    # longer referenced by these cells.
    $rng.Style = "Normal"
    $rng.Value = $g.Value
}
